$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 3 through 17 (old data rows beyond the new single data row)
$ws.Range("A3:B17").EntireRow.Delete() | Out-Null

# Update the remaining data row (row 2) with new values
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = 2.061751933828537
